$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows (Chiang Rai x3, Songkhla 17-25, Sisaket, Chonburi) from bottom to top
# so earlier row indices remain valid while deleting.
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(2).Delete()
